$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 262, shifting existing rows 262-324 down to 263-325
$ws.Rows.Item(262).Insert()

# Populate the newly inserted row 262 with its data
$ws.Cells.Item(262, 1).Value = 5
$ws.Cells.Item(262, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(262, 3).Value = "Maule"
$ws.Cells.Item(262, 4).Value = 44754
$ws.Cells.Item(262, 5).Value = 7
$ws.Cells.Item(262, 6).Value = 100112003
$ws.Cells.Item(262, 7).Value = "Ajo"
$ws.Cells.Item(262, 8).Value = "Chino"
$ws.Cells.Item(262, 9).Value = "Primera"
$ws.Cells.Item(262, 10).Value = 200
$ws.Cells.Item(262, 11).Value = 19000
$ws.Cells.Item(262, 12).Value = 19000
$ws.Cells.Item(262, 13).Value = 19000
$ws.Cells.Item(262, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(262, 15).Value = "China"
$ws.Cells.Item(262, 16).Value = 1900
$ws.Cells.Item(262, 17).Value = 10
$ws.Cells.Item(262, 18).Value = "Hortaliza"
